# Generate Report for Handoff
# Update status text "In Translation" -> "Ready for handoff" and refresh the
# handoff timestamps on the Overview, zh-cn and de-de sheets, then widen the
# affected status/date columns to fit the new text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status text: "In Translation" -> "Ready for handoff" -------------------
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value     = "Ready for handoff"
$wsDeDe.Range("C2").Value     = "Ready for handoff"

# --- Timestamps ---------------------------------------------------------
# Overview "Latest HO Xliff Generate Date" and de-de "Latest Handoff Datetime"
$wsOverview.Range("G2").Value = "2016-08-30 12:48:25"
$wsDeDe.Range("H2").Value     = "2016-08-30 12:48:25"

# zh-cn "Latest Handoff Datetime"
$wsZhCn.Range("H2").Value = "2016-08-30 12:48:19"

# --- Column widths (grow to fit "Ready for handoff") ---------------------
$wsOverview.Columns.Item(5).ColumnWidth = 17.2159881591797
$wsOverview.Columns.Item(6).ColumnWidth = 17.2159881591797
$wsZhCn.Columns.Item(3).ColumnWidth     = 17.2159881591797
$wsDeDe.Columns.Item(3).ColumnWidth     = 17.2159881591797
